# Update odds values on Sheet1 to reflect the latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("O2").Value = 1.67
$ws.Range("P2").Value = 2.1

# Row 4
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6

# Row 5
$ws.Range("J5").Value = 2.2
$ws.Range("K5").Value = 2.25
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.9
$ws.Range("U5").Value = 1.95
$ws.Range("V5").Value = 1.8
$ws.Range("X5").Value = 7.5
$ws.Range("AB5").Value = 26
$ws.Range("AC5").Value = 10
$ws.Range("AE5").Value = 17
$ws.Range("AF5").Value = 51
$ws.Range("AM5").Value = 301
$ws.Range("AN5").Value = 3.6
$ws.Range("BB5").Value = 251

# Row 6
$ws.Range("Q6").Value = 1.93
$ws.Range("R6").Value = 1.93

# Row 7
$ws.Range("G7").Value = 2.15
$ws.Range("H7").Value = 2.9
$ws.Range("J7").Value = 2.7
$ws.Range("L7").Value = 4.1
$ws.Range("S7").Value = 1.44
$ws.Range("T7").Value = 2.42
$ws.Range("W7").Value = 6.6
$ws.Range("X7").Value = 10
$ws.Range("Z7").Value = 21
$ws.Range("AA7").Value = 18.5
$ws.Range("AG7").Value = 8.5
$ws.Range("AH7").Value = 18
$ws.Range("AI7").Value = 12.5
$ws.Range("AK7").Value = 40
$ws.Range("AL7").Value = 50
$ws.Range("AM7").Value = 800
$ws.Range("AO7").Value = 11
$ws.Range("AP7").Value = 18.5
$ws.Range("AR7").Value = 70
$ws.Range("AU7").Value = 6.8
$ws.Range("AX7").Value = 21
$ws.Range("AY7").Value = 28
$ws.Range("AZ7").Value = 110
$ws.Range("BB7").Value = 350

# Row 8
$ws.Range("G8").Value = 2.45
$ws.Range("I8").Value = 2.8
$ws.Range("J8").Value = 3.2
$ws.Range("L8").Value = 3.6
$ws.Range("S8").Value = 1.47
$ws.Range("T8").Value = 2.5
$ws.Range("X8").Value = 11
$ws.Range("Z8").Value = 23
$ws.Range("AK8").Value = 26
$ws.Range("AR8").Value = 67
$ws.Range("AT8").Value = 2.5
$ws.Range("AY8").Value = 29

# Row 9
$ws.Range("S9").Value = 1.47

# Row 10
$ws.Range("G10").Value = 2.4
$ws.Range("H10").Value = 2.7

# Row 11
$ws.Range("L11").Value = 3.75
$ws.Range("X11").Value = 11
$ws.Range("AD11").Value = 6
$ws.Range("AH11").Value = 17
$ws.Range("AI11").Value = 12
$ws.Range("AJ11").Value = 34
$ws.Range("AK11").Value = 26
$ws.Range("AN11").Value = 4.33
$ws.Range("AO11").Value = 12
$ws.Range("AZ11").Value = 51
$ws.Range("BB11").Value = 151

# Row 12
$ws.Range("G12").Value = 1.29
$ws.Range("H12").Value = 5
$ws.Range("I12").Value = 11
$ws.Range("J12").Value = 1.8
$ws.Range("K12").Value = 2.38
$ws.Range("L12").Value = 10
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 9.5
$ws.Range("Q12").Value = 1.95
$ws.Range("R12").Value = 1.9
$ws.Range("U12").Value = 2.5
$ws.Range("V12").Value = 1.5
$ws.Range("X12").Value = 5.5
$ws.Range("Y12").Value = 9.5
$ws.Range("Z12").Value = 7.5
$ws.Range("AD12").Value = 10
$ws.Range("AE12").Value = 29
$ws.Range("AF12").Value = 101
$ws.Range("AG12").Value = 19
$ws.Range("AI12").Value = 29
$ws.Range("AK12").Value = 81
$ws.Range("AL12").Value = 81
$ws.Range("AP12").Value = 21
$ws.Range("AQ12").Value = 17
$ws.Range("AU12").Value = 11
$ws.Range("AV12").Value = 81
$ws.Range("AZ12").Value = 301
